# Adding more performance numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (F2/G2, the two "Cannot be run" shared-string cells) is dropped
# entirely. Clear() removes the now-empty <row> element (and the save
# path drops the now-unreferenced shared strings automatically).
$ws.Range("F2:G2").Clear()

# New "Average run time (25 Users)" figures for the rows that previously
# had no value in column F.
$ws.Range("F3").Value = 17107
$ws.Range("F4").Value = 7385
$ws.Range("F5").Value = 8422
$ws.Range("F7").Value = 24064
$ws.Range("F8").Value = 36568

# Update the active selection/scroll position shown in the sheet view.
$ws.Range("F38").Select()
